$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "Alex"
$ws.Range("C9").Value = "Alex"
$ws.Range("C5").Value = "Fait"
$ws.Range("C4").Value = "Vicky"
$ws.Range("C3").Value = "Nico"

$ws.Range("C4").Select()
